{"js": "// Replace the 25 \"three-digit \u00f7 one-digit\" answer cells in the practice\n// table with the new values from the commit. Every old value is unique in\n// the document, so body.search() + Range.insertText(..., \"Replace\") is a\n// safe, order-independent way to apply every substitution.\nconst replacements = [\n  [\"348\u00f77=49, 5\", \"218\u00f77=31, 1\"],\n  [\"958\u00f78=119, 6\", \"540\u00f73=180, 0\"],\n  [\"948\u00f74=237, 0\", \"165\u00f74=41, 1\"],\n  [\"300\u00f72=150, 0\", \"249\u00f79=27, 6\"],\n  [\"675\u00f73=225, 0\", \"804\u00f79=89, 3\"],\n  [\"914\u00f74=228, 2\", \"985\u00f78=123, 1\"],\n  [\"740\u00f77=105, 5\", \"344\u00f79=38, 2\"],\n  [\"617\u00f75=123, 2\", \"221\u00f72=110, 1\"],\n  [\"769\u00f79=85, 4\", \"134\u00f72=67, 0\"],\n  [\"114\u00f72=57, 0\", \"192\u00f79=21, 3\"],\n  [\"511\u00f73=170, 1\", \"456\u00f73=152, 0\"],\n  [\"453\u00f79=50, 3\", \"910\u00f75=182, 0\"],\n  [\"956\u00f79=106, 2\", \"543\u00f78=67, 7\"],\n  [\"748\u00f75=149, 3\", \"286\u00f75=57, 1\"],\n  [\"462\u00f73=154, 0\", \"576\u00f77=82, 2\"],\n  [\"450\u00f77=64, 2\", \"655\u00f74=163, 3\"],\n  [\"358\u00f73=119, 1\", \"480\u00f77=68, 4\"],\n  [\"590\u00f76=98, 2\", \"939\u00f77=134, 1\"],\n  [\"872\u00f77=124, 4\", \"863\u00f77=123, 2\"],\n  [\"615\u00f74=153, 3\", \"270\u00f72=135, 0\"],\n  [\"499\u00f72=249, 1\", \"172\u00f79=19, 1\"],\n  [\"502\u00f77=71, 5\", \"913\u00f75=182, 3\"],\n  [\"949\u00f72=474, 1\", \"909\u00f77=129, 6\"],\n  [\"362\u00f77=51, 5\", \"675\u00f77=96, 3\"],\n  [\"664\u00f72=332, 0\", \"128\u00f79=14, 2\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Could not find text to replace: \"${oldText}\"`);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the 25 \"three-digit \u00f7 one-digit\" answer cells in the practice\n# table with the new values from the commit. Every old value is unique in\n# the document, so a simple Find/Replace (wdReplaceAll semantics, but each\n# string only ever matches once) safely applies every substitution without\n# depending on row/column indices.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"348\u00f77=49, 5\", \"218\u00f77=31, 1\"),\n    @(\"958\u00f78=119, 6\", \"540\u00f73=180, 0\"),\n    @(\"948\u00f74=237, 0\", \"165\u00f74=41, 1\"),\n    @(\"300\u00f72=150, 0\", \"249\u00f79=27, 6\"),\n    @(\"675\u00f73=225, 0\", \"804\u00f79=89, 3\"),\n    @(\"914\u00f74=228, 2\", \"985\u00f78=123, 1\"),\n    @(\"740\u00f77=105, 5\", \"344\u00f79=38, 2\"),\n    @(\"617\u00f75=123, 2\", \"221\u00f72=110, 1\"),\n    @(\"769\u00f79=85, 4\", \"134\u00f72=67, 0\"),\n    @(\"114\u00f72=57, 0\", \"192\u00f79=21, 3\"),\n    @(\"511\u00f73=170, 1\", \"456\u00f73=152, 0\"),\n    @(\"453\u00f79=50, 3\", \"910\u00f75=182, 0\"),\n    @(\"956\u00f79=106, 2\", \"543\u00f78=67, 7\"),\n    @(\"748\u00f75=149, 3\", \"286\u00f75=57, 1\"),\n    @(\"462\u00f73=154, 0\", \"576\u00f77=82, 2\"),\n    @(\"450\u00f77=64, 2\", \"655\u00f74=163, 3\"),\n    @(\"358\u00f73=119, 1\", \"480\u00f77=68, 4\"),\n    @(\"590\u00f76=98, 2\", \"939\u00f77=134, 1\"),\n    @(\"872\u00f77=124, 4\", \"863\u00f77=123, 2\"),\n    @(\"615\u00f74=153, 3\", \"270\u00f72=135, 0\"),\n    @(\"499\u00f72=249, 1\", \"172\u00f79=19, 1\"),\n    @(\"502\u00f77=71, 5\", \"913\u00f75=182, 3\"),\n    @(\"949\u00f72=474, 1\", \"909\u00f77=129, 6\"),\n    @(\"362\u00f77=51, 5\", \"675\u00f77=96, 3\"),\n    @(\"664\u00f72=332, 0\", \"128\u00f79=14, 2\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $rng = $d.Content\n    $found = $rng.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n\n    if (-not $found) {\n        throw \"Could not find text to replace: $oldText\"\n    }\n}\n"}
